## Generate Report for Handback
## -----------------------------------------------------------------------
## This script mirrors the localization-status report being regenerated
## after a handback: the "Ready for handoff" status becomes "Handed back:
## in sync with en-US", the per-language sheets gain "Latest Target File"
## / "Latest Handback DateTime" data (with a hyperlink on the target file
## name, same as the source-file hyperlink in column A), and a handful of
## columns are widened to fit the new text.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/80ff9cdf4fbdc1cf93d4287f22288f62ce3c9593/e2e/"

# ---------------------------------------------------------------------
# 1. Overview sheet: status text + column widths
# ---------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")

$ovw.Range("E2").Value = "Handed back: in sync with en-US"
$ovw.Range("F2").Value = "Handed back: in sync with en-US"
$ovw.Range("E3").Value = "Handed back: in sync with en-US"
$ovw.Range("F3").Value = "Handed back: in sync with en-US"

$ovw.Columns.Item(5).ColumnWidth = 29.166666666666668
$ovw.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# 2. zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "Handed back: in sync with en-US"

$zh.Range("I2").Value = "a.md"
$zh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-31 04:38:45"

$zh.Range("I3").Value = "a.md"
$zh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-31 04:38:45"

$zh.Hyperlinks.Add($zh.Range("I2"), ($repoBase + "a.md"), "", "", "a.md")
$zh.Hyperlinks.Add($zh.Range("I3"), ($repoBase + "a.md"), "", "", "a.md")

$zh.Columns.Item(3).ColumnWidth = 29.166666666666668
$zh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# 3. de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "Handed back: in sync with en-US"

$de.Range("I2").Value = "a.md"
$de.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$de.Range("K2").Value = "2016-08-31 04:38:52"

$de.Range("I3").Value = "a.md"
$de.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$de.Range("K3").Value = "2016-08-31 04:38:52"

$de.Hyperlinks.Add($de.Range("I2"), ($repoBase + "a.md"), "", "", "a.md")
$de.Hyperlinks.Add($de.Range("I3"), ($repoBase + "a.md"), "", "", "a.md")

$de.Columns.Item(3).ColumnWidth = 29.166666666666668
$de.Columns.Item(10).ColumnWidth = 39.166666666666664
